$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(39, 8).Value = 660.61536
$ws.Cells.Item(39, 9).Value = 189.81818
$ws.Cells.Item(39, 10).Value = 3250
$ws.Cells.Item(39, 11).Value = 569.4545400000001
$ws.Cells.Item(39, 12).Value = 9750
$ws.Cells.Item(39, 13).Value = -273.4545400000001
$ws.Cells.Item(39, 14).Value = -10342

$ws.Cells.Item(40, 8).Value = 836923.8
$ws.Cells.Item(40, 9).Value = 1254103.6
$ws.Cells.Item(40, 10).Value = 2564.25
$ws.Cells.Item(40, 11).Value = 1254103.6
$ws.Cells.Item(40, 12).Value = 2564.25
$ws.Cells.Item(40, 13).Value = -1253928.6
$ws.Cells.Item(40, 14).Value = -2914.25

$ws.Cells.Item(53, 8).Value = 559.7727
$ws.Cells.Item(53, 9).Value = 307.25
$ws.Cells.Item(53, 11).Value = 307.25
$ws.Cells.Item(53, 13).Value = 329.75

$ws.Cells.Item(112, 8).Value = 2314.5386
$ws.Cells.Item(112, 9).Value = 966.1539
$ws.Cells.Item(112, 10).Value = 3662.923
$ws.Cells.Item(112, 11).Value = 2898.4617
$ws.Cells.Item(112, 12).Value = 10988.769
$ws.Cells.Item(112, 13).Value = -1790.4617
$ws.Cells.Item(112, 14).Value = -13204.769

$ws.Cells.Item(118, 8).Value = 347.33334
$ws.Cells.Item(118, 9).Value = 265.75
$ws.Cells.Item(118, 11).Value = 797.25
$ws.Cells.Item(118, 13).Value = 859.75

$ws.Cells.Item(132, 8).Value = 5111.0264
$ws.Cells.Item(132, 9).Value = 4980.543
$ws.Cells.Item(132, 11).Value = 14941.629
$ws.Cells.Item(132, 13).Value = -12411.629

$ws.Cells.Item(136, 8).Value = 130780
$ws.Cells.Item(136, 10).Value = 130780
$ws.Cells.Item(136, 12).Value = 130780
$ws.Cells.Item(136, 14).Value = -140980

$ws.Cells.Item(137, 8).Value = 5048.7407
$ws.Cells.Item(137, 9).Value = 2569.5715
$ws.Cells.Item(137, 11).Value = 7708.7145
$ws.Cells.Item(137, 13).Value = -5158.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1858.4149
$ws.Cells.Item(32, 9).Value = 997.6512
$ws.Cells.Item(32, 11).Value = 997.6512
$ws.Cells.Item(32, 13).Value = -710.6512

$ws.Cells.Item(132, 8).Value = 37101084
$ws.Cells.Item(132, 9).Value = 13595.708
$ws.Cells.Item(132, 11).Value = 40787.124
$ws.Cells.Item(132, 13).Value = -38257.124

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(25, 8).Value = 288.15384
$ws.Cells.Item(25, 9).Value = 278.91666
$ws.Cells.Item(25, 10).Value = 399
$ws.Cells.Item(25, 11).Value = 278.91666
$ws.Cells.Item(25, 12).Value = 399
$ws.Cells.Item(25, 13).Value = -43.91665999999998
$ws.Cells.Item(25, 14).Value = -869

$ws.Cells.Item(99, 8).Value = 3530
$ws.Cells.Item(99, 9).Value = 2794.2222
$ws.Cells.Item(99, 11).Value = 2794.2222
$ws.Cells.Item(99, 13).Value = -1296.2222

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 3302.4
$ws.Cells.Item(16, 9).Value = 3499.6667
$ws.Cells.Item(16, 11).Value = 3499.6667
$ws.Cells.Item(16, 13).Value = -3212.6667

$ws.Cells.Item(31, 8).Value = 5685508.5
$ws.Cells.Item(31, 9).Value = 2305.3635
$ws.Cells.Item(31, 10).Value = 11368712
$ws.Cells.Item(31, 11).Value = 2305.3635
$ws.Cells.Item(31, 12).Value = 11368712
$ws.Cells.Item(31, 13).Value = -2010.3635
$ws.Cells.Item(31, 14).Value = -11369302

$ws.Cells.Item(34, 8).Value = 5685508.5
$ws.Cells.Item(34, 9).Value = 2305.3635
$ws.Cells.Item(34, 10).Value = 11368712
$ws.Cells.Item(34, 11).Value = 2305.3635
$ws.Cells.Item(34, 12).Value = 11368712
$ws.Cells.Item(34, 13).Value = -2103.3635
$ws.Cells.Item(34, 14).Value = -11369116

$ws.Cells.Item(113, 8).Value = 3302.4
$ws.Cells.Item(113, 9).Value = 3499.6667
$ws.Cells.Item(113, 11).Value = 3499.6667
$ws.Cells.Item(113, 13).Value = -1329.6667

$ws.Cells.Item(122, 8).Value = 8773354
$ws.Cells.Item(122, 9).Value = 1648.3334
$ws.Cells.Item(122, 11).Value = 4945.0002
$ws.Cells.Item(122, 13).Value = -2495.0002

$ws.Cells.Item(132, 8).Value = 71147.62
$ws.Cells.Item(132, 9).Value = 89012.35000000001
$ws.Cells.Item(132, 10).Value = 2666.1667
$ws.Cells.Item(132, 11).Value = 267037.05
$ws.Cells.Item(132, 12).Value = 7998.500100000001
$ws.Cells.Item(132, 13).Value = -264507.05
$ws.Cells.Item(132, 14).Value = -13058.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(109, 8).Value = 6122.909
$ws.Cells.Item(109, 9).Value = 10875
$ws.Cells.Item(109, 11).Value = 32625
$ws.Cells.Item(109, 13).Value = -31585

$ws.Cells.Item(124, 8).Value = 4126.2
$ws.Cells.Item(124, 9).Value = 2815
$ws.Cells.Item(124, 11).Value = 8445
$ws.Cells.Item(124, 13).Value = -3535

$ws.Cells.Item(131, 8).Value = 54009.773
$ws.Cells.Item(131, 9).Value = 87560.336
$ws.Cells.Item(131, 10).Value = 13749.1
$ws.Cells.Item(131, 11).Value = 262681.008
$ws.Cells.Item(131, 12).Value = 41247.3
$ws.Cells.Item(131, 13).Value = -257641.008
$ws.Cells.Item(131, 14).Value = -51327.3

$ws.Cells.Item(132, 8).Value = 3351.7046
$ws.Cells.Item(132, 9).Value = 2306.2188
$ws.Cells.Item(132, 11).Value = 20755.9692
$ws.Cells.Item(132, 13).Value = -18225.9692

$ws.Cells.Item(134, 8).Value = 5254.1113
$ws.Cells.Item(134, 9).Value = 1041.4286
$ws.Cells.Item(134, 11).Value = 3124.2858
$ws.Cells.Item(134, 13).Value = 1945.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 1091980.8
$ws.Cells.Item(2, 9).Value = 2953020.8
$ws.Cells.Item(2, 10).Value = 1026.138
$ws.Cells.Item(2, 11).Value = 2953020.8
$ws.Cells.Item(2, 12).Value = 1026.138
$ws.Cells.Item(2, 13).Value = -2952907.8
$ws.Cells.Item(2, 14).Value = -1252.138

$ws.Cells.Item(62, 8).Value = 79331.664
$ws.Cells.Item(62, 10).Value = 78995
$ws.Cells.Item(62, 12).Value = 78995
$ws.Cells.Item(62, 14).Value = -80367

$ws.Cells.Item(65, 8).Value = 79331.664
$ws.Cells.Item(65, 10).Value = 78995
$ws.Cells.Item(65, 12).Value = 236985
$ws.Cells.Item(65, 14).Value = -243849

$ws.Cells.Item(70, 8).Value = 85697.42
$ws.Cells.Item(70, 9).Value = 108386.7
$ws.Cells.Item(70, 11).Value = 108386.7
$ws.Cells.Item(70, 13).Value = -108116.7

$ws.Cells.Item(73, 8).Value = 85697.42
$ws.Cells.Item(73, 9).Value = 108386.7
$ws.Cells.Item(73, 11).Value = 108386.7
$ws.Cells.Item(73, 13).Value = -107450.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(20, 8).Value = 20000
$ws.Cells.Item(20, 9).Value = 20000
$ws.Cells.Item(20, 11).Value = 20000
$ws.Cells.Item(20, 13).Value = -19774

$ws.Cells.Item(40, 8).Value = 2771.2778
$ws.Cells.Item(40, 9).Value = 2696.697
$ws.Cells.Item(40, 11).Value = 2696.697
$ws.Cells.Item(40, 13).Value = -2560.697

$ws.Cells.Item(46, 8).Value = 2036.6177
$ws.Cells.Item(46, 9).Value = 989.8
$ws.Cells.Item(46, 11).Value = 989.8
$ws.Cells.Item(46, 13).Value = -801.8

$ws.Cells.Item(55, 8).Value = 689.13635
$ws.Cells.Item(55, 9).Value = 238.36363
$ws.Cells.Item(55, 10).Value = 1139.909
$ws.Cells.Item(55, 11).Value = 238.36363
$ws.Cells.Item(55, 12).Value = 1139.909
$ws.Cells.Item(55, 13).Value = -65.36363
$ws.Cells.Item(55, 14).Value = -1485.909

$ws.Cells.Item(61, 8).Value = 2788.4167
$ws.Cells.Item(61, 10).Value = 5406.3335
$ws.Cells.Item(61, 12).Value = 5406.3335
$ws.Cells.Item(61, 14).Value = -5810.3335

$ws.Cells.Item(93, 8).Value = 362704.8
$ws.Cells.Item(93, 9).Value = 1484.1714
$ws.Cells.Item(93, 10).Value = 2168808
$ws.Cells.Item(93, 11).Value = 1484.1714
$ws.Cells.Item(93, 12).Value = 2168808
$ws.Cells.Item(93, 13).Value = -236.1713999999999
$ws.Cells.Item(93, 14).Value = -2171304

$ws.Cells.Item(113, 8).Value = 2788.4167
$ws.Cells.Item(113, 10).Value = 5406.3335
$ws.Cells.Item(113, 12).Value = 5406.3335
$ws.Cells.Item(113, 14).Value = -9746.333500000001

$ws.Cells.Item(132, 8).Value = 5725
$ws.Cells.Item(132, 9).Value = 5346.3335
$ws.Cells.Item(132, 11).Value = 16039.0005
$ws.Cells.Item(132, 13).Value = -13509.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(40, 8).Value = 14642.143
$ws.Cells.Item(40, 9).Value = 7500
$ws.Cells.Item(40, 10).Value = 19998.75
$ws.Cells.Item(40, 11).Value = 7500
$ws.Cells.Item(40, 12).Value = 19998.75
$ws.Cells.Item(40, 13).Value = -7351
$ws.Cells.Item(40, 14).Value = -20296.75

$ws.Cells.Item(122, 8).Value = 11115615
$ws.Cells.Item(122, 9).Value = 3315.6667
$ws.Cells.Item(122, 11).Value = 9947.000100000001
$ws.Cells.Item(122, 13).Value = -7497.000100000001

$ws.Cells.Item(132, 8).Value = 2398.4285
$ws.Cells.Item(132, 9).Value = 1920.3334
$ws.Cells.Item(132, 10).Value = 3832.7144
$ws.Cells.Item(132, 11).Value = 5761.0002
$ws.Cells.Item(132, 12).Value = 11498.1432
$ws.Cells.Item(132, 13).Value = -3231.0002
$ws.Cells.Item(132, 14).Value = -16558.1432

$ws.Cells.Item(135, 8).Value = 16747524
$ws.Cells.Item(135, 10).Value = 16747524
$ws.Cells.Item(135, 12).Value = 16747524
$ws.Cells.Item(135, 14).Value = -16757664
